# 自动更新Excel文件 - 2025-11-23 23:11:59
# Decrement the "剩余" (days remaining, column E) counter for every shop
# row by one day. When a row's counter has bottomed out (reaches 1, i.e.
# the last remaining day), the delivery cycle restarts: the counter is
# reset to a fresh 7-day cycle and the "开始时间" (start date, column F)
# is advanced by 7 days to the new cycle's start. Rows whose start date
# is not a well-formed 8-digit yyyyMMdd value (data-entry errors) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
if (-not $lastRow) { $lastRow = 99 }

for ($row = 2; $row -le $lastRow; $row++) {
    $remainingCell = $ws.Cells.Item($row, 5)   # column E - 剩余
    $startCell = $ws.Cells.Item($row, 6)       # column F - 开始时间

    $remaining = $remainingCell.Value()
    $startDate = $startCell.Value()

    if ($remaining -eq $null -or $startDate -eq $null) {
        continue
    }

    # Skip malformed start dates (not a plain 8-digit yyyyMMdd number).
    $startText = [string]$startDate
    if ($startText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $remainingCell.Value = 7
        $startCell.Value = $startDate + 7
    } else {
        $remainingCell.Value = $remaining - 1
    }
}
